$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - "Change in inventories"
$ws.Range("B6").Value = -147300000.0
$ws.Range("C6").Value = -114500000.0
$ws.Range("D6").Value = -80300000.0
$ws.Range("E6").Value = -47400000.0
$ws.Range("F6").Value = -40900000.0
$ws.Range("G6").Value = -49100000.0

# Row 7 - "Change in payables and accrued liability"
$ws.Range("B7").Value = 416000000.0
$ws.Range("C7").Value = 529000000.0
$ws.Range("D7").Value = 468000000.0
$ws.Range("E7").Value = 341500000.0
$ws.Range("F7").Value = 284600000.0
$ws.Range("G7").Value = 179000000.0

# Row 28 - "Capital Stock Change": B28 goes from empty inline string to numeric value
$ws.Range("B28").Value = 17300000.0
